$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data for columns B (Circulation), C (ILL Loans), D (ILL Borrows)
# for rows 3 through 59, matching the library order in column A.
$data = @(
    @(49589,7209,8406),
    @(25960,2961,2962),
    @(86096,7410,7672),
    @(1883,732,180),
    @(54308,8996,7099),
    @(5936,1521,1223),
    @(6700,1296,772),
    @(3150,450,299),
    @(402,241,29),
    @(3,0,0),
    @(1097,255,307),
    @(3452,1390,1046),
    @(5732,2151,923),
    @(3937,1496,633),
    @(2549,986,218),
    @(19548,3062,3616),
    @(1606,681,479),
    @(21472,2650,3575),
    @(211,426,20),
    @(19386,2616,3405),
    @(1245,459,175),
    @(22472,3016,4183),
    @(88283,8370,10606),
    @(6835,2180,1073),
    @(0,0,0),
    @(6137,1267,1426),
    @(1641,492,326),
    @(16927,2965,2916),
    @(552,187,240),
    @(3015,1961,433),
    @(18865,3601,3138),
    @(11843,3374,2440),
    @(6640,705,1570),
    @(65546,6221,6403),
    @(9586,3140,1418),
    @(28175,2218,3246),
    @(1190,1062,206),
    @(2310,511,813),
    @(2734,328,101),
    @(9805,557,273),
    @(299,114,67),
    @(810,58,66),
    @(1045,14,7),
    @(3744,991,470),
    @(13870,3922,2370),
    @(35146,3821,4804),
    @(17207,4042,1393),
    @(12716,1341,1806),
    @(36309,3408,5185),
    @(5593,668,1393),
    @(15728,3362,2620),
    @(2231,1458,838),
    @(2304,1430,168),
    @(4215,1073,1208),
    @(14218,5292,2909),
    @(14693,1071,514),
    @(768035,114845,106245)
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
}
